$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.568.32"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.35%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.35"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.83%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.59"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.17%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.11%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4798"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.06%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2810"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.88%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06498"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.90%  "

# Row 10 - WrappedEther
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.977.89"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.47%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07466"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.83%  "

# Row 12 - Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.48"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.93%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.084"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.32%  "

# Row 14 - Litecoin
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.83"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.11%  "

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6624"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.17%  "

# Row 16 - WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.522.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.29%  "

# Row 17 - Avalanche
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.22"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.48%  "

# Row 18 - Dai
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.09%  "

# Row 19 - swap ShibaInu -> WrappedliquidstakedEther2.0
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.223.72"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.17%  "

# Row 20 - swap WrappedliquidstakedEther2.0 -> ShibaInu
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007542"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.22%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "227.13"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.32%  "

# Row 22 - BinanceUSD
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.03%  "

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.258"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.43%  "

# Row 24 - Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.145"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.40%  "

# Row 25 - Cosmos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.292"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.32%  "

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.09"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.34%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.41"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.23%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.938"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.55%  "

# Row 29 - Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.399"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.69%  "

# Row 30 - Stellar
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09656"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.18%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.328"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.23%  "

# Row 32 - Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.992"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.15%  "

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05040"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.15%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  +4.41%  "

# Row 35 - ImmutableX
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7477"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.16%  "

# Row 36 - HuobiToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.713"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.73%  "

# Row 37 - VeChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01853"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.16%  "

# Row 38 - MXToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.636"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.39%  "

# Row 39 - TrustWalletToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9101"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.80%  "

# Row 40 - RenderToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.065"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.16%  "

# Row 41 - Quant
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.99"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.52%  "

# Row 42 - TheSandbox
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4258"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.12%  "

# Row 43 - swap PaxDollar -> FraxShare
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.757"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.31%  "

# Row 44 - swap FraxShare -> PaxDollar
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9989"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.48%  "

# Row 45 - Aptos
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.338"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.05%  "

# Row 46 - Algorand
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1281"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.70%  "

# Row 47 - Aave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "63.87"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.31%  "

# Row 48 - NEARProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.467"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.88%  "

# Row 49 - EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.913"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.59%  "

# Row 50 - Elrond
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.60"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.55%  "

# Row 51 - Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05654"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.32%  "
